# Update the built-in Heading styles so every heading level shares one
# uniform accent color, and adjust the weight/size/emphasis of a few
# levels, per the "Change styles in reference.docx" commit:
#
#   - Heading 1: drop the 181 (0xB5) theme shade on its accent1 color so
#                it matches the plain accent1 color used everywhere else.
#   - Heading 2: 16 pt -> 14 pt.
#   - Heading 3: 14 pt -> 12 pt.
#   - Heading 4: bold  -> italic.
#   - Heading 5: no longer forces italic explicitly (Heading 4 now owns
#                that look); complex-script italic flag is left alone.

$d = $word.ActiveDocument

# wdThemeColorAccent1, as used by the other heading styles in this
# document's color object model.
$wdThemeColorAccent1 = 4

# --- Heading 1: strip the 181 (0xB5) theme shade -----------------------
$h1 = $d.Styles("Heading1")
$h1.Font.TextColor.ObjectThemeColor = $wdThemeColorAccent1

# --- Heading 2: 16 pt -> 14 pt ------------------------------------------
$h2 = $d.Styles("Heading2")
$h2.Font.Size = 14
$h2.Font.SizeBi = 14

# --- Heading 3: 14 pt -> 12 pt ------------------------------------------
$h3 = $d.Styles("Heading3")
$h3.Font.Size = 12
$h3.Font.SizeBi = 12

# --- Heading 4: bold -> italic -------------------------------------------
$h4 = $d.Styles("Heading4")
$h4.Font.Bold = $False
$h4.Font.Italic = $True

# --- Heading 5: drop the explicit italic flag ---------------------------
$h5 = $d.Styles("Heading5")
$h5.Font.Italic = $False
